$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style representing the default (unstyled) data-cell format,
# used to strip the "quote prefix" style Excel applies when a numeric-
# looking literal is entered with a leading apostrophe (forces Text type
# without Excel silently re-parsing it back into a Number).
$defaultStyle = $ws.Range("D20").Style

$ws.Range("D2").Formula = "'285.10"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Formula = "'2.98%"
$ws.Range("E2").Style = $defaultStyle
$ws.Range("D3").Formula = "'28.73"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Formula = "'5.96%"
$ws.Range("E3").Style = $defaultStyle
$ws.Range("D4").Formula = "'4.918"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Formula = "'1.29%"
$ws.Range("E4").Style = $defaultStyle
$ws.Range("D5").Formula = "'0.06493"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Formula = "'1.41%"
$ws.Range("E5").Style = $defaultStyle
$ws.Range("D6").Formula = "'7.217"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Formula = "'4.00%"
$ws.Range("E6").Style = $defaultStyle
$ws.Range("D7").Formula = "'1.331"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Formula = "'10.38%"
$ws.Range("E7").Style = $defaultStyle
$ws.Range("D8").Formula = "'0.9116"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Formula = "'3.99%"
$ws.Range("E8").Style = $defaultStyle
$ws.Range("E9").Formula = "'0.73%"
$ws.Range("E9").Style = $defaultStyle
$ws.Range("D10").Formula = "'0.06431"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Formula = "'25.49%"
$ws.Range("E10").Style = $defaultStyle
$ws.Range("D11").Formula = "'0.07626"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Formula = "'1.60%"
$ws.Range("E11").Style = $defaultStyle
$ws.Range("D12").Formula = "'0.02985"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Formula = "'0.81%"
$ws.Range("E12").Style = $defaultStyle
$ws.Range("D13").Formula = "'0.08968"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Formula = "'-0.15%"
$ws.Range("E13").Style = $defaultStyle
$ws.Range("D14").Formula = "'0.001599"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Formula = "'1.73%"
$ws.Range("E14").Style = $defaultStyle
$ws.Range("D15").Formula = "'0.0006531"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Formula = "'2.41%"
$ws.Range("E15").Style = $defaultStyle
$ws.Range("D16").Formula = "'0.006022"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Formula = "'-1.22%"
$ws.Range("E16").Style = $defaultStyle
$ws.Range("E17").Formula = "'-0.47%"
$ws.Range("E17").Style = $defaultStyle
$ws.Range("D18").Formula = "'3.371"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Formula = "'1.99%"
$ws.Range("E18").Style = $defaultStyle
$ws.Range("E19").Formula = "'-1.46%"
$ws.Range("E19").Style = $defaultStyle
$ws.Range("E21").Formula = "'0.22%"
$ws.Range("E21").Style = $defaultStyle
$ws.Range("D22").Formula = "'3.977"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Formula = "'1.90%"
$ws.Range("E22").Style = $defaultStyle
$ws.Range("E23").Formula = "'12.73%"
$ws.Range("E23").Style = $defaultStyle
$ws.Range("D24").Formula = "'0.04475"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Formula = "'1.10%"
$ws.Range("E24").Style = $defaultStyle
$ws.Range("D25").Formula = "'0.001186"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Formula = "'0.87%"
$ws.Range("E25").Style = $defaultStyle
$ws.Range("D26").Formula = "'0.004325"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Formula = "'12.00%"
$ws.Range("E26").Style = $defaultStyle
$ws.Range("E28").Formula = "'-9.14%"
$ws.Range("E28").Style = $defaultStyle
$ws.Range("D29").Formula = "'0.0001637"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Formula = "'-15.66%"
$ws.Range("E29").Style = $defaultStyle
$ws.Range("D40").Formula = "'0.04166"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Formula = "'0.26%"
$ws.Range("E40").Style = $defaultStyle
$ws.Range("D41").Formula = "'0.006731"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Formula = "'-1.04%"
$ws.Range("E41").Style = $defaultStyle
$ws.Range("D42").Formula = "'0.1234"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Formula = "'5.19%"
$ws.Range("E42").Style = $defaultStyle
$ws.Range("D43").Formula = "'0.002142"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Formula = "'9.32%"
$ws.Range("E43").Style = $defaultStyle
$ws.Range("D44").Formula = "'0.01177"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Formula = "'-1.21%"
$ws.Range("E44").Style = $defaultStyle
$ws.Range("D45").Formula = "'0.00005434"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Formula = "'2.55%"
$ws.Range("E45").Style = $defaultStyle
$ws.Range("E46").Formula = "'20.98%"
$ws.Range("E46").Style = $defaultStyle
$ws.Range("D47").Formula = "'0.01852"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Formula = "'0.00%"
$ws.Range("E47").Style = $defaultStyle
